$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H15 ("innards in microwave for a short moment(s)") — fix typo/wording: drop trailing "s"
$ws.Range("N18").Value = "innards in microwave for a short moment"

# H7 ("Big spider crawls away after sight") — tiny wording improvement: "spider(s)"
$ws.Range("N10").Value = "Big spider(s) crawls away after sight"

# H13 previously had a placeholder idea ("in game spider jumps directly to the camera").
# Reuse/relocate the stronger idea that used to live on H20 and give it a base completion
# value plus a supporting detail note (the "added base for H13" part of the commit).
$ws.Range("N16").Value = "Big Spider jumps through gamePlayers face into metaPlayers face"
$ws.Range("Q16").Value = 0.7
$ws.Range("R16").Value = "lower immersedvalue before so you aren't too close to the screen, fine tuning, reaction, crawling away"

# H20 no longer needs its own separate description now that it lives on H13; clear it but
# keep the percentage-style formatting on its %DoneSat cell (empty, still formatted as 0%).
$ws.Range("N23").ClearContents()
$ws.Range("Q23").NumberFormat = "0%"

# Leave the cursor where the author's last edit landed.
$null = $ws.Range("N9").Select()
